# Re-generate statistics with fixed minutes/seconds formatting (zero-padded to
# two digits) in the "haul" (Квала) column of the top-by-best-speed sheet.
#
# Time strings look like: "<H> ч. <M> мин. <S> сек." where H, M, S are
# numbers without leading zeros in the original export. The fix pads the
# minutes and seconds components to two digits, e.g.
#   "11 ч. 22 мин. 9 сек."  -> "11 ч. 22 мин. 09 сек."
#   "1 ч. 5 мин. 27 сек."   -> "1 ч. 05 мин. 27 сек."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$haulColumn = 9   # column I = "Квала" (haul time)
$changed = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $haulColumn)
    $text = $cell.Text

    if ([string]::IsNullOrEmpty($text)) {
        continue
    }

    if ($text -match '^(\d+) ч\. (\d+) мин\. (\d+) сек\.$') {
        $hours = $matches[1]
        $minutes = $matches[2]
        $seconds = $matches[3]

        $paddedMinutes = $minutes.PadLeft(2, '0')
        $paddedSeconds = $seconds.PadLeft(2, '0')

        $newText = "$hours ч. $paddedMinutes мин. $paddedSeconds сек."

        if ($newText -ne $text) {
            $cell.Value = $newText
            $changed++
        }
    }
}

Write-Host "Updated $changed haul time cells."
